# Update cryptos price/volume snapshot (GitHub Actions scheduled refresh).
# Each row's Price (column D) and Volume(1h) (column E) are refreshed with
# the latest scraped text. Numeric-looking price strings are written with a
# leading quote (quote-prefix) so Excel keeps them as text, matching the
# original inline-string cells (prices use '.' both as thousands separator
# and decimal point, e.g. '25.915.29', so they must stay text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '25.915.29'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.742.85'
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").Value = '''0.9998'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''247.00'
$ws.Range("E5").Value = '  +4.75%  '
$ws.Range("D6").Value = '''1.0000'
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = '''0.5036'
$ws.Range("E7").Value = '  -4.72%  '
$ws.Range("D8").Value = '''0.2726'
$ws.Range("E8").Value = '  -2.99%  '
$ws.Range("D9").Value = '''0.06183'
$ws.Range("E9").Value = '  -0.26%  '
$ws.Range("D10").Value = '1.751.72'
$ws.Range("E10").Value = '  +0.29%  '
$ws.Range("D11").Value = '''0.07240'
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").Value = '''0.6510'
$ws.Range("E12").Value = '  +0.49%  '
$ws.Range("D13").Value = '''15.13'
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("D15").Value = '''77.42'
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = '''1.001'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '''0.9996'
$ws.Range("E17").Value = '  +0.03%  '
$ws.Range("D18").Value = '25.941.74'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Value = '''0.000006799'
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D21").Value = '1.980.18'
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").Value = '''4.318'
$ws.Range("E22").Value = '  -0.47%  '
$ws.Range("D23").Value = '''8.654'
$ws.Range("E23").Value = '  -1.14%  '
$ws.Range("D24").Value = '''5.397'
$ws.Range("E24").Value = '  +2.74%  '
$ws.Range("D25").Value = '''136.46'
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").Value = '''1.504'
$ws.Range("E26").Value = '  -1.06%  '
$ws.Range("D27").Value = '''15.19'
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("D28").Value = '''1.767'
$ws.Range("E28").Value = '  -2.81%  '
$ws.Range("D29").Value = '''105.60'
$ws.Range("E29").Value = '  +0.81%  '
$ws.Range("D30").Value = '''3.919'
$ws.Range("E30").Value = '  +2.70%  '
$ws.Range("D31").Value = '''0.08234'
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").Value = '''3.630'
$ws.Range("E32").Value = '  -0.92%  '
$ws.Range("D33").Value = '''0.04674'
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").Value = '''2.656'
$ws.Range("D35").Value = '''0.9932'
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").Value = '''0.6177'
$ws.Range("E36").Value = '  -3.07%  '
$ws.Range("D37").Value = '''2.737'
$ws.Range("E37").Value = '  +1.06%  '
$ws.Range("D38").Value = '''0.01604'
$ws.Range("E38").Value = '  -0.30%  '
$ws.Range("D39").Value = '''1.915'
$ws.Range("E39").Value = '  -3.45%  '
$ws.Range("D40").Value = '''0.9996'
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").Value = '''99.46'
$ws.Range("E41").Value = '  -2.21%  '
$ws.Range("D42").Value = '''0.3873'
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("D44").Value = '''4.989'
$ws.Range("E44").Value = '  -0.87%  '
$ws.Range("D45").Value = '''0.1141'
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("D46").Value = '''6.288'
$ws.Range("E46").Value = '  -1.64%  '
$ws.Range("D47").Value = '''55.49'
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("D48").Value = '''0.05237'
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("D49").Value = '''30.56'
$ws.Range("E49").Value = '  -1.32%  '
$ws.Range("D50").Value = '''7.541'
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("D51").Value = '''0.3415'
$ws.Range("E51").Value = '  -1.82%  '
